$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '89.517.43'
$ws.Range('E2').Value = '  +1.25%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.042.96'
$ws.Range('E3').Value = '  -2.52%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.28'
$ws.Range('E5').Value = '  -1.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '611.57'
$ws.Range('E6').Value = '  -3.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.363'
$ws.Range('E7').Value = '  -7.48%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.881'
$ws.Range('E8').Value = '  +20.76%  '
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.042.07'
$ws.Range('E10').Value = '  -2.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.663'
$ws.Range('E11').Value = '  +19.93%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.186'
$ws.Range('E12').Value = '  +4.30%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000237'
$ws.Range('E13').Value = '  -5.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.37'
$ws.Range('E14').Value = '  +1.81%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '88.389.50'
$ws.Range('E15').Value = '  +0.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '31.89'
$ws.Range('E16').Value = '  -0.92%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.600.68'
$ws.Range('E17').Value = '  -2.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.033.92'
$ws.Range('E18').Value = '  -2.67%  '
$ws.Range('E19').Value = '  +1.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000211'
$ws.Range('E20').Value = '  -2.98%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.40'
$ws.Range('E21').Value = '  +1.80%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '425.30'
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.00'
$ws.Range('E23').Value = '  +2.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.14'
$ws.Range('E24').Value = '  -2.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.41'
$ws.Range('E25').Value = '  +1.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '83.44'
$ws.Range('E26').Value = '  +4.90%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.71'
$ws.Range('E27').Value = '  +2.24%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.198.69'
$ws.Range('E28').Value = '  -2.20%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('B30').Value = 'Cronos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.163'
$ws.Range('E30').Value = '  +2.65%  '
$ws.Range('B31').Value = 'Binance-PegBSC-USD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.01'
$ws.Range('E31').Value = '  +2.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.20'
$ws.Range('E32').Value = '  +0.61%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '503.37'
$ws.Range('E33').Value = '  -1.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.60'
$ws.Range('E34').Value = '  -7.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.63'
$ws.Range('E35').Value = '  -4.30%  '
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '22.77'
$ws.Range('E36').Value = '  +4.49%  '
$ws.Range('B37').Value = 'PancakeSwap'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.79'
$ws.Range('E37').Value = '  -1.95%  '
$ws.Range('E38').Value = '  -4.03%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.130'
$ws.Range('E39').Value = '  +3.99%  '
$ws.Range('B40').Value = 'WhiteBITCoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '22.26'
$ws.Range('E40').Value = '  +0.33%  '
$ws.Range('E41').Value = '  +0.14%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('B43').Value = 'PolygonEcosystemToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.362'
$ws.Range('E43').Value = '  +0.21%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.137'
$ws.Range('E44').Value = '  +8.90%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.82'
$ws.Range('E45').Value = '  -1.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '147.07'
$ws.Range('E46').Value = '  +2.34%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '43.36'
$ws.Range('E47').Value = '  -0.59%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0679'
$ws.Range('E48').Value = '  +11.33%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.07'
$ws.Range('E49').Value = '  +3.92%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.22'
$ws.Range('E50').Value = '  +3.33%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '158.04'
$ws.Range('E51').Value = '  -3.92%  '
